$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '69.337.46'
$ws.Range("E2").Value = '  +0.05%  '

$ws.Range("D3").Value = '3.672.57'
$ws.Range("E3").Value = '  -0.45%  '

$ws.Range("E4").Value = '  -0.02%  '

$ws.Range("D5").Value = "'684.66"
$ws.Range("E5").Value = '  +0.05%  '

$ws.Range("D6").Value = "'158.83"
$ws.Range("E6").Value = '  -2.24%  '

$ws.Range("E7").Value = '  -0.03%  '

$ws.Range("E8").Value = '  -1.32%  '

$ws.Range("E9").Value = '  -2.28%  '

$ws.Range("D10").Value = "'7.03"
$ws.Range("E10").Value = '  -2.63%  '

$ws.Range("D11").Value = "'0.435"
$ws.Range("E11").Value = '  -3.48%  '

$ws.Range("E12").Value = '  -1.92%  '

$ws.Range("D13").Value = '4.292.59'
$ws.Range("E13").Value = '  -0.45%  '

$ws.Range("D14").Value = "'32.18"
$ws.Range("E14").Value = '  -3.95%  '

$ws.Range("B15").Value = 'WrappedEther'
$ws.Range("C15").Value = 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
$ws.Range("D15").Value = '3.693.95'
$ws.Range("E15").Value = '  +0.14%  '

$ws.Range("B16").Value = 'WrappedBTC'
$ws.Range("C16").Value = 'https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc'
$ws.Range("D16").Value = '69.332.10'
$ws.Range("E16").Value = '  -0.09%  '

$ws.Range("E17").Value = '  +2.08%  '

$ws.Range("E18").Value = '  -3.22%  '

$ws.Range("D19").Value = "'6.38"
$ws.Range("E19").Value = '  -3.61%  '

$ws.Range("D20").Value = "'469.00"
$ws.Range("E20").Value = '  -2.81%  '

$ws.Range("E21").Value = '  +1.49%  '

$ws.Range("D22").Value = "'0.648"
$ws.Range("E22").Value = '  -2.59%  '

$ws.Range("D23").Value = "'79.69"
$ws.Range("E23").Value = '  -0.36%  '

$ws.Range("D24").Value = '3.818.39'
$ws.Range("E24").Value = '  -0.41%  '

$ws.Range("E25").Value = '  +0.01%  '

$ws.Range("D26").Value = "'0.0000123"
$ws.Range("E26").Value = '  -4.60%  '

$ws.Range("D27").Value = "'10.88"
$ws.Range("E27").Value = '  -5.37%  '

$ws.Range("E28").Value = '  -4.10%  '

$ws.Range("E29").Value = '  -1.99%  '

$ws.Range("E30").Value = '  -6.20%  '

$ws.Range("B31").Value = 'Binance-PegBSC-USD'
$ws.Range("C31").Value = 'https://coinranking.com/coin/i5jggxiwp+binance-pegbsc-usd-bsc-usd'
$ws.Range("D31").Value = "'1.00"
$ws.Range("E31").Value = '  +0.10%  '

$ws.Range("B32").Value = 'NEARProtocol'
$ws.Range("C32").Value = 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
$ws.Range("D32").Value = "'6.54"
$ws.Range("E32").Value = '  -3.67%  '

$ws.Range("E33").Value = '  -5.91%  '

$ws.Range("D34").Value = "'26.84"
$ws.Range("E34").Value = '  -0.65%  '

$ws.Range("D35").Value = '3.646.13'
$ws.Range("E35").Value = '  -0.18%  '

$ws.Range("D36").Value = "'0.160"
$ws.Range("E36").Value = '  -3.16%  '

$ws.Range("E37").Value = '  -4.85%  '

$ws.Range("D38").Value = "'6.09"
$ws.Range("E38").Value = '  +0.84%  '

$ws.Range("E40").Value = '  +1.94%  '

$ws.Range("E41").Value = '  -4.76%  '

$ws.Range("E42").Value = '  -0.07%  '

$ws.Range("D43").Value = "'166.58"
$ws.Range("E43").Value = '  +5.82%  '

$ws.Range("E44").Value = '  -2.17%  '

$ws.Range("D45").Value = "'47.49"
$ws.Range("E45").Value = '  -1.42%  '

$ws.Range("D46").Value = "'0.000282"
$ws.Range("E46").Value = '  +1.44%  '

$ws.Range("B47").Value = 'dogwifhat'
$ws.Range("C47").Value = 'https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif'
$ws.Range("D47").Value = "'2.70"
$ws.Range("E47").Value = '  -3.94%  '

$ws.Range("B48").Value = 'SuiNetwork'
$ws.Range("C48").Value = 'https://coinranking.com/coin/3xJluUMvp+suinetwork-sui'
$ws.Range("D48").Value = "'1.11"
$ws.Range("E48").Value = '  +4.06%  '

$ws.Range("E49").Value = '  -0.28%  '

$ws.Range("D50").Value = "'27.28"
$ws.Range("E50").Value = '  -2.28%  '

$ws.Range("D51").Value = "'7.77"
$ws.Range("E51").Value = '  -3.71%  '
